$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update measurement rows 2-5 with the new data set ---
$row2 = @(45162.50694444445,4.347,5.032,0.777,10.022,8.942,4.483,9.151999999999999,5.852,2.498,6.03,7.398,4.257,0.73,3.371,5.736,3.045,0.138,0.227,54.637,11.069,2.959,6.58,4.913,0.755,5.641,3.12,5.177,2.517,6.419,0.192,7.159,2.015,4.399)
for ($i = 0; $i -lt $row2.Length; $i++) { $ws.Cells.Item(2, $i + 1).Value = $row2[$i] }

$row3 = @(45162.51388888889,21.689,17.023,1.054,47.359,39.476,17.47,60.885,26.593,11.944,18.721,20.641,20.119,5.243,16.855,24.717,14.233,0.178,0.654,255.532,48.194,15.447,32.299,17.684,2.332,31.455,13.93,13.415,14.161,21.005,0.115,55.147,9.172000000000001,19.72)
for ($i = 0; $i -lt $row3.Length; $i++) { $ws.Cells.Item(3, $i + 1).Value = $row3[$i] }

$row4 = @(45162.52083333334,14.625,11.448,0.712,31.93,26.604,11.755,47.072,17.952,8.082000000000001,12.588,13.893,13.558,3.533,11.35,16.687,9.613,0.111,0.417,169.846,32.601,10.388,21.869,11.948,1.559,23.166,9.375999999999999,9.021000000000001,9.561999999999999,14.136,0.092,42.787,6.226,13.275)
for ($i = 0; $i -lt $row4.Length; $i++) { $ws.Cells.Item(4, $i + 1).Value = $row4[$i] }

$row5 = @(45162.52777777778,6.04,4.88,0.37,13.24,11.09,4.93,22.86,7.51,3.41,5.43,6.12,5.62,1.41,4.66,7.04,4.03,0.04,0.15,66.92,13.79,4.22,9.220000000000001,5.16,0.66,11.01,3.89,4.01,3.93,6.08,0.08,20.89,2.66,5.5)
for ($i = 0; $i -lt $row5.Length; $i++) { $ws.Cells.Item(5, $i + 1).Value = $row5[$i] }

# --- Drop the trailing measurement row (old row 6); dimension shrinks to A1:AH5 ---
$ws.Rows.Item(6).Delete()

# --- Widen a handful of data columns from 7 to 8 characters ---
# (ColumnWidth is offset from the stored column width; -0.85 lands exactly on
# the target integer width after Excel's internal pixel rounding.)
$wideColumns = @(3, 10, 15, 17, 22, 28, 29)
foreach ($c in $wideColumns) {
    $ws.Columns.Item($c).ColumnWidth = 7.15
}
